$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching style of existing headers (bold, centered, bordered)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

$ws.Cells.Item(1, 8).Copy()
$ws.Range($ws.Cells.Item(1, 9), $ws.Cells.Item(1, 10)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row number, I value, J value
$data = @(
    @(2,7,7),
    @(3,8,8),
    @(4,8,8),
    @(5,8,8),
    @(6,7,8),
    @(7,8,8),
    @(8,7,7),
    @(9,8,8),
    @(10,8,8),
    @(11,7,7),
    @(12,9,9),
    @(13,9,9),
    @(14,9,9),
    @(15,8,9),
    @(16,9,9),
    @(17,9,9),
    @(18,9,9),
    @(19,9,9),
    @(20,9,9),
    @(21,8,8),
    @(22,8,8),
    @(23,7,7),
    @(24,7,7),
    @(25,7,8),
    @(26,8,8),
    @(27,8,8),
    @(28,8,8),
    @(29,7,8),
    @(30,8,8),
    @(31,7,7),
    @(32,8,8),
    @(33,8,8),
    @(34,9,9),
    @(35,9,9),
    @(36,8,9),
    @(37,7,8),
    @(38,6,6),
    @(39,7,8),
    @(40,6,6),
    @(41,8,8),
    @(42,5,5),
    @(43,3,4),
    @(44,4,5),
    @(45,6,6),
    @(46,5,5),
    @(47,5,5),
    @(48,7,7),
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
